$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91..161 down to 92..162 (copy values bottom-up to avoid overwrite)
for ($r = 162; $r -ge 92; $r--) {
    $src = $r - 1
    $ws.Range("C$r").Value2 = $ws.Range("C$src").Value2
    $ws.Range("D$r").Value2 = $ws.Range("D$src").Value2
    $ws.Range("E$r").Value2 = $ws.Range("E$src").Value2
    $ws.Range("F$r").Value2 = $ws.Range("F$src").Value2
    $ws.Range("G$r").Value2 = $ws.Range("G$src").Value2
}

# New row 91: NTLite (write alias "ntl" first so it lands before "NTLite" in the shared-string table)
$ws.Range("D91").Value2 = "ntl"
$ws.Range("C91").Value2 = "NTLite"
$ws.Range("E91").Value2 = "[Application]"
$ws.Range("F91").Value2 = $null
$ws.Range("G91").Value2 = $null

# Extend the row-numbering dynamic array formula to cover the new row (B3:B161 -> B3:B162)
$ws.Range("B3").Formula2 = "=SEQUENCE(COUNTA(C:C)-1)"

# Keep a trailing blank row after the table (row 166 mirrors the old trailing blank row 165),
# extending the used range down by one row to match the shift.
$ws.Range("B165:G165").Copy()
$ws.Range("B166:G166").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "done"
